$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 650.1539
$ws.Range("J2").Value = 975.7143
$ws.Range("L2").Value = 975.7143
$ws.Range("N2").Value = -1201.7143
$ws.Range("H17").Value = 1327.3549
$ws.Range("J17").Value = 1327.3549
$ws.Range("L17").Value = 3982.0647
$ws.Range("N17").Value = -4318.0647
$ws.Range("H33").Value = 192.6923
$ws.Range("I33").Value = 198.72728
$ws.Range("J33").Value = 159.5
$ws.Range("K33").Value = 198.72728
$ws.Range("L33").Value = 159.5
$ws.Range("M33").Value = 30.27271999999999
$ws.Range("N33").Value = -617.5
$ws.Range("H58").Value = 1248.4286
$ws.Range("I58").Value = 119.5
$ws.Range("J58").Value = 1700
$ws.Range("K58").Value = 358.5
$ws.Range("L58").Value = 5100
$ws.Range("M58").Value = -208.5
$ws.Range("N58").Value = -5400
$ws.Range("H74").Value = 5946.394
$ws.Range("J74").Value = 6328.6313
$ws.Range("L74").Value = 6328.6313
$ws.Range("N74").Value = -8200.631300000001
$ws.Range("H77").Value = 5946.394
$ws.Range("J77").Value = 6328.6313
$ws.Range("L77").Value = 31643.1565
$ws.Range("N77").Value = -41003.1565
$ws.Range("H88").Value = 1460.3125
$ws.Range("J88").Value = 1655.4166
$ws.Range("L88").Value = 1655.4166
$ws.Range("N88").Value = -2467.4166
$ws.Range("H91").Value = 1460.3125
$ws.Range("J91").Value = 1655.4166
$ws.Range("L91").Value = 1655.4166
$ws.Range("N91").Value = -4463.4166
$ws.Range("H93").Value = 66999
$ws.Range("J93").Value = 66999
$ws.Range("L93").Value = 66999
$ws.Range("N93").Value = -71991
$ws.Range("H106").Value = 4790122
$ws.Range("I106").Value = 6499384.5
$ws.Range("K106").Value = 6499384.5
$ws.Range("M106").Value = -6498753.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 2922.3
$ws.Range("I46").Value = 4166.6665
$ws.Range("J46").Value = 2389
$ws.Range("K46").Value = 4166.6665
$ws.Range("L46").Value = 2389
$ws.Range("M46").Value = -3847.6665
$ws.Range("N46").Value = -3027
$ws.Range("H61").Value = 3699.7334
$ws.Range("I61").Value = 2277
$ws.Range("K61").Value = 2277
$ws.Range("M61").Value = -2065
$ws.Range("H74").Value = 2851.5173
$ws.Range("I74").Value = 2811.423
$ws.Range("K74").Value = 2811.423
$ws.Range("M74").Value = -1937.423
$ws.Range("H77").Value = 2851.5173
$ws.Range("I77").Value = 2811.423
$ws.Range("K77").Value = 14057.115
$ws.Range("M77").Value = -9689.114999999998
$ws.Range("H136").Value = 3699.7334
$ws.Range("I136").Value = 2277
$ws.Range("K136").Value = 6831
$ws.Range("M136").Value = -4281

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1872.0769
$ws.Range("I20").Value = 1747.375
$ws.Range("J20").Value = 2071.6
$ws.Range("K20").Value = 1747.375
$ws.Range("L20").Value = 2071.6
$ws.Range("M20").Value = -1500.375
$ws.Range("N20").Value = -2565.6
$ws.Range("H86").Value = 6638.2563
$ws.Range("I86").Value = 3873.037
$ws.Range("K86").Value = 3873.037
$ws.Range("M86").Value = -2750.037
$ws.Range("H89").Value = 6638.2563
$ws.Range("I89").Value = 3873.037
$ws.Range("K89").Value = 19365.185
$ws.Range("M89").Value = -13749.185
$ws.Range("H106").Value = 30129.625
$ws.Range("J106").Value = 30129.625
$ws.Range("L106").Value = 30129.625
$ws.Range("N106").Value = -32653.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2572.261
$ws.Range("I31").Value = 2907.1428
$ws.Range("J31").Value = 2425.75
$ws.Range("K31").Value = 2907.1428
$ws.Range("L31").Value = 2425.75
$ws.Range("M31").Value = -2612.1428
$ws.Range("N31").Value = -3015.75
$ws.Range("H34").Value = 2572.261
$ws.Range("I34").Value = 2907.1428
$ws.Range("J34").Value = 2425.75
$ws.Range("K34").Value = 2907.1428
$ws.Range("L34").Value = 2425.75
$ws.Range("M34").Value = -2705.1428
$ws.Range("N34").Value = -2829.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 80.94118
$ws.Range("I2").Value = 83.30768999999999
$ws.Range("J2").Value = 73.25
$ws.Range("K2").Value = 499.84614
$ws.Range("L2").Value = 439.5
$ws.Range("M2").Value = -386.84614
$ws.Range("N2").Value = -665.5
$ws.Range("H29").Value = 133.25
$ws.Range("I29").Value = 219.5
$ws.Range("K29").Value = 658.5
$ws.Range("M29").Value = -381.5
$ws.Range("H37").Value = 141933.88
$ws.Range("J37").Value = 141933.88
$ws.Range("L37").Value = 425801.64
$ws.Range("N37").Value = -426025.64
$ws.Range("H68").Value = 2509.6
$ws.Range("I68").Value = 2183.3333
$ws.Range("K68").Value = 6549.999899999999
$ws.Range("M68").Value = -5738.999899999999
$ws.Range("H71").Value = 2509.6
$ws.Range("I71").Value = 2183.3333
$ws.Range("K71").Value = 19649.9997
$ws.Range("M71").Value = -15593.9997
$ws.Range("H76").Value = 11277.444
$ws.Range("I76").Value = 6832.6665
$ws.Range("K76").Value = 20497.9995
$ws.Range("M76").Value = -20114.9995
$ws.Range("H79").Value = 11277.444
$ws.Range("I79").Value = 6832.6665
$ws.Range("K79").Value = 20497.9995
$ws.Range("M79").Value = -19171.9995
$ws.Range("H81").Value = 17666.5
$ws.Range("I81").Value = 16199.8
$ws.Range("J81").Value = 25000
$ws.Range("K81").Value = 48599.39999999999
$ws.Range("L81").Value = 75000
$ws.Range("M81").Value = -47476.39999999999
$ws.Range("N81").Value = -77246
$ws.Range("H84").Value = 17666.5
$ws.Range("I84").Value = 16199.8
$ws.Range("J84").Value = 25000
$ws.Range("K84").Value = 145798.2
$ws.Range("L84").Value = 225000
$ws.Range("M84").Value = -140182.2
$ws.Range("N84").Value = -236232
$ws.Range("H131").Value = 23457970
$ws.Range("J131").Value = 37175652
$ws.Range("L131").Value = 111526956
$ws.Range("N131").Value = -111537036
$ws.Range("H140").Value = 5328428
$ws.Range("I140").Value = 11365325
$ws.Range("K140").Value = 34095975
$ws.Range("M140").Value = -34090795

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 54047.5
$ws.Range("J98").Value = 54047.5
$ws.Range("L98").Value = 54047.5
$ws.Range("N98").Value = -60037.5
$ws.Range("H119").Value = 99800
$ws.Range("J119").Value = 99800
$ws.Range("L119").Value = 99800
$ws.Range("N119").Value = -109476
$ws.Range("H128").Value = 60000
$ws.Range("J128").Value = 60000
$ws.Range("L128").Value = 60000
$ws.Range("N128").Value = -69960
$ws.Range("H141").Value = 115321.22
$ws.Range("I141").Value = 95000
$ws.Range("J141").Value = 121127.29
$ws.Range("K141").Value = 95000
$ws.Range("L141").Value = 121127.29
$ws.Range("M141").Value = -89820
$ws.Range("N141").Value = -131487.29

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1935.5714
$ws.Range("I46").Value = 1414.4286
$ws.Range("J46").Value = 2456.7144
$ws.Range("K46").Value = 1414.4286
$ws.Range("L46").Value = 2456.7144
$ws.Range("M46").Value = -1226.4286
$ws.Range("N46").Value = -2832.7144
$ws.Range("H104").Value = 30498
$ws.Range("J104").Value = 30498
$ws.Range("L104").Value = 30498
$ws.Range("N104").Value = -37486
$ws.Range("H131").Value = 73950
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 3041.875
$ws.Range("I132").Value = 2491.9688
$ws.Range("J132").Value = 3775.0833
$ws.Range("K132").Value = 7475.9064
$ws.Range("L132").Value = 11325.2499
$ws.Range("M132").Value = -4945.9064
$ws.Range("N132").Value = -16385.2499
$ws.Range("H140").Value = 83326.71000000001
$ws.Range("J140").Value = 90482.836
$ws.Range("L140").Value = 90482.836
$ws.Range("N140").Value = -100842.836

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 55648
$ws.Range("J46").Value = 55648
$ws.Range("L46").Value = 55648
$ws.Range("N46").Value = -56110
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("H123").Value = 131998.6
$ws.Range("J123").Value = 131998.6
$ws.Range("L123").Value = 131998.6
$ws.Range("N123").Value = -141798.6
$ws.Range("H124").Value = 58937.25
$ws.Range("I124").Value = 50000
$ws.Range("J124").Value = 61916.332
$ws.Range("K124").Value = 50000
$ws.Range("L124").Value = 61916.332
$ws.Range("M124").Value = -45090
$ws.Range("N124").Value = -71736.33199999999
$ws.Range("H125").Value = 131048.4
$ws.Range("J125").Value = 131048.4
$ws.Range("L125").Value = 131048.4
$ws.Range("N125").Value = -140888.4
$ws.Range("H134").Value = 55648
$ws.Range("J134").Value = 55648
$ws.Range("L134").Value = 166944
$ws.Range("N134").Value = -172014
$ws.Range("H135").Value = 35662.816
$ws.Range("J135").Value = 35662.816
$ws.Range("L135").Value = 35662.816
$ws.Range("N135").Value = -45802.816
$ws.Range("H136").Value = 7953.7812
$ws.Range("J136").Value = 24235
$ws.Range("L136").Value = 72705
$ws.Range("N136").Value = -77805
$ws.Range("H140").Value = 91475
$ws.Range("J140").Value = 91475
$ws.Range("L140").Value = 91475
$ws.Range("N140").Value = -101835
$ws.Range("H141").Value = 77111.94
$ws.Range("J141").Value = 77111.94
$ws.Range("L141").Value = 77111.94
$ws.Range("N141").Value = -87471.94
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
